$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date value for every existing data row
# (rows 2..51) from 2023-09-21 (45190) to 2023-09-23 (45192).
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}

# Row 51 picks up an explicit row height once the sheet grows past it.
$ws.Rows.Item(51).RowHeight = 15

# Append the new record as row 52.
$ws.Cells.Item(52, 1).Value = "A 44646-2023"

$ws.Cells.Item(52, 2).Value = 45189
$ws.Cells.Item(52, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(52, 3).Value = 45192
$ws.Cells.Item(52, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(52, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item(52, 5).Value = "LUND"

$ws.Cells.Item(52, 7).Value = 3.2
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 13).Value = 0
$ws.Cells.Item(52, 14).Value = 0
$ws.Cells.Item(52, 15).Value = 0
$ws.Cells.Item(52, 16).Value = 0
$ws.Cells.Item(52, 17).Value = 0

$ws.Cells.Item(52, 18).Value = ""
$ws.Cells.Item(52, 18).WrapText = $true
